# Date and Time pickers enabled.
# Add a new row (#11) to the task list documenting that the user can now
# change the Date and Time the Task is due.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 11
$ws.Range("C13").Value = "Allow the user to change the Date and Time the Task is Due"

# Match the author's final selection (cursor left on the newly added cell).
[void]$ws.Range("C13").Select()
